# bug fixes in count pipeline and updated counting data
#
# 1) "main" sheet: recompute `area` (column C) for a set of rows, and for
#    the rows whose area shrank enough to no longer straddle two boxes,
#    drop `weighted_ct` (column D) from 2 to 1.
# 2) "conjoined" sheet: the rows corresponding to those now-unconjoined
#    ids are removed (the sheet only lists frame/id pairs that are still
#    conjoined).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("main")

# row -> (new area, new weighted_ct-or-$null)
$updates = @(
    @(18,   16120, $null),
    @(22,   53130, 1),
    @(110,  52614, $null),
    @(167,  46494, 1),
    @(282,  68094, $null),
    @(300,  21140, $null),
    @(347,  82348, $null),
    @(353,  71370, $null),
    @(401,  63174, $null),
    @(416,  42768, 1),
    @(446,  34456, 1),
    @(460,  22270, $null),
    @(609,  22040, $null),
    @(619,  19995, $null),
    @(777,  32035, $null),
    @(781,  53732, $null),
    @(810,  58706, 1),
    @(816,  41328, 1),
    @(889,  19028, $null),
    @(892,  32025, $null),
    @(940,  30240, $null),
    @(1021, 72518, $null),
    @(1022, 43433, 1),
    @(1047, 72927, $null),
    @(1225, 71604, $null),
    @(1259, 35696, $null),
    @(1307, 38793, 1),
    @(1324, 58045, 1),
    @(1339, 28542, $null),
    @(1437, 28416, $null)
)

foreach ($u in $updates) {
    $row = $u[0]
    $area = $u[1]
    $wct = $u[2]
    $ws.Cells.Item($row, 3).Value = $area
    if ($null -ne $wct) {
        $ws.Cells.Item($row, 4).Value = $wct
    }
}

# "conjoined" sheet: remove the rows for the ids that are no longer
# conjoined (weighted_ct dropped 2 -> 1 above). Delete from the bottom
# up so earlier row numbers stay valid as we go.
$wsConj = $wb.Worksheets.Item("conjoined")

$rowsToDelete = @(73,71,70,68,65,64,57,55,54,52,51,50,49,46,45,44,43,35,32,29,26,25,24,21,20,19,18,17,10,7,4,3)

foreach ($r in $rowsToDelete) {
    $wsConj.Rows.Item($r).Delete()
}
